# Generate Report for Handback
# Updates the "zh-cn" and "de-de" report sheets for row 7
# (750cf41e-eb23-4615-a393-16bdff136d60) with the handback results:
#  - Latest Target File (I) becomes a hyperlink to the .md file that was
#    actually handed back
#  - Latest Handback File (J) is filled in with the generated xlf name
#  - Latest Handback DateTime (K) is filled in
#  - Error Detail (P) reports that the handback file version is stale
# Also widens the "Error Detail" column (P) so the message is readable.

$wb = $excel.ActiveWorkbook

$currentUrl = 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7f4e3a7d196487b90da8e764f92d67f14330d789/e2e/750cf41e-eb23-4615-a393-16bdff136d60.md'
$errorDetail = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7f4e3a7d196487b90da8e764f92d67f14330d789/e2e/750cf41e-eb23-4615-a393-16bdff136d60.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/33fe81a380562a3c29ac80ab30804d5013d48ef1/e2e/750cf41e-eb23-4615-a393-16bdff136d60.md.'

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $currentUrl, "", "", "750cf41e-eb23-4615-a393-16bdff136d60.md")
$wsZh.Range("J7").Value = "750cf41e-eb23-4615-a393-16bdff136d60.f133367a032665177a688607fd4d32d62012a1cd.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-31 16:49:45"
$wsZh.Range("P7").Value = $errorDetail

$wsZh.Columns.Item(16).ColumnWidth = 40 - 5/6

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $currentUrl, "", "", "750cf41e-eb23-4615-a393-16bdff136d60.md")
$wsDe.Range("J7").Value = "750cf41e-eb23-4615-a393-16bdff136d60.f133367a032665177a688607fd4d32d62012a1cd.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-31 16:49:52"
$wsDe.Range("P7").Value = $errorDetail

$wsDe.Columns.Item(16).ColumnWidth = 40 - 5/6
